$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Species types" sheet: add a new row for water (H2O) as a metabolite.
# ---------------------------------------------------------------------
$wsSpeciesTypes = $wb.Worksheets.Item("Species types")

$wsSpeciesTypes.Range("A8").Value = "H2O"
$wsSpeciesTypes.Range("A8").Style = "Normal"
$wsSpeciesTypes.Range("A8").HorizontalAlignment = -4131
$wsSpeciesTypes.Range("A8").VerticalAlignment = -4160
$wsSpeciesTypes.Range("A8").WrapText = $true

$wsSpeciesTypes.Range("B8").Value = "H2O"
$wsSpeciesTypes.Range("B8").Style = "Normal"
$wsSpeciesTypes.Range("B8").HorizontalAlignment = -4131
$wsSpeciesTypes.Range("B8").VerticalAlignment = -4160
$wsSpeciesTypes.Range("B8").WrapText = $true

$wsSpeciesTypes.Range("C8").Value = "ASP_test_2016_2"
$wsSpeciesTypes.Range("C8").Style = "Normal"

$wsSpeciesTypes.Range("F8").Value = 18.0152
$wsSpeciesTypes.Range("F8").Style = "Normal"
$wsSpeciesTypes.Range("F8").HorizontalAlignment = -4131
$wsSpeciesTypes.Range("F8").VerticalAlignment = -4160
$wsSpeciesTypes.Range("F8").WrapText = $true

$wsSpeciesTypes.Range("H8").Value = "metabolite"

$wsSpeciesTypes.Activate()
[void]$wsSpeciesTypes.Range("A8").Select()

# ---------------------------------------------------------------------
# "Concentrations" sheet: add initial concentrations for H2O in the
# extracellular space and cytoplasm compartments.
# ---------------------------------------------------------------------
$wsConcentrations = $wb.Worksheets.Item("Concentrations")

$wsConcentrations.Range("A8").Value = "H2O[e]"
$wsConcentrations.Range("A8").Style = "Normal"
$wsConcentrations.Range("A8").HorizontalAlignment = -4131
$wsConcentrations.Range("A8").VerticalAlignment = -4160
$wsConcentrations.Range("A8").WrapText = $true
$wsConcentrations.Range("B8").Value = 1

$wsConcentrations.Range("A9").Value = "H2O[c]"
$wsConcentrations.Range("A9").Style = "Normal"
$wsConcentrations.Range("A9").HorizontalAlignment = -4131
$wsConcentrations.Range("A9").VerticalAlignment = -4160
$wsConcentrations.Range("A9").WrapText = $true
$wsConcentrations.Range("B9").Value = 1

$wsConcentrations.Activate()
[void]$wsConcentrations.Range("B10").Select()

# ---------------------------------------------------------------------
# "Reactions" sheet keeps its own zoom level; it is no longer the
# active tab once Concentrations becomes selected above.
# ---------------------------------------------------------------------
$wsReactions = $wb.Worksheets.Item("Reactions")
$wsReactions.Activate()
$excel.ActiveWindow.Zoom = 120

# Re-activate Concentrations so it ends up as the active tab.
$wsConcentrations.Activate()
